# Update "想去人数" (want-to-go count) figures for the 苏州-漫展信息 workbook.
# These edits mirror a data refresh pulled from the source site; the same
# rows appear on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($name -eq "展览") {
        $ws.Range("F5").Value  = 13398
        $ws.Range("F13").Value = 14517
        $ws.Range("F27").Value = 5557
        $ws.Range("F30").Value = 5356
        $ws.Range("F32").Value = 22
        $ws.Range("F33").Value = 140
    }
    else {
        $ws.Range("F5").Value  = 13398
        $ws.Range("F14").Value = 14517
        $ws.Range("F28").Value = 5557
        $ws.Range("F31").Value = 5356
        $ws.Range("F33").Value = 22
        $ws.Range("F34").Value = 140
    }
}
